$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new (leftmost) date columns, pushing the existing
#     B:E columns (and everything they hold) three positions to the
#     right, landing on E:H. This mirrors the weekly "roll forward"
#     of the MarketBeat watch report where the newest dates are
#     always inserted right after the analyst name column.
$ws.Columns("B:D").Insert()

# --- Header row: the three freshly inserted columns get the new
#     date labels.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Fill the new columns for every existing analyst row with the
#     "UN" (unchanged) placeholder used throughout the sheet.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- Argus (row 9) picked up a new initiation on both of the newly
#     added mid-week columns.
$ws.Range("C9").Value = "6/25/2018,Initiates,Buy,`$275.00"
$ws.Range("D9").Value = "6/25/2018,Initiates,Buy,`$275.00"

# --- Two brand new analysts/firms added to the bottom of the watch
#     list this week.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

$wb.Save()
